# Update leave date after add test result
# - Remove the extra quarantine-day rows (rows 3-8), keeping only the
#   single remaining record in row 2.
# - Update the remaining record: ID 165 -> 57, test date 2021-04-05 (44291)
#   -> 2021-12-10 (44540).
# - Flip the test result for that record from "am tinh" (negative) to
#   "duong tinh" (positive).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3 through 8 entirely (shifts nothing below them up since they
# are the last rows of data).
$ws.Rows("3:8").Delete()

# Update the surviving data row (row 2).
$ws.Range("A2").Value = 57
$ws.Range("C2").Value = 44540

# Change this row's result text last, once it is the only cell left that
# references the old shared string, so the string itself is updated in
# place rather than creating a duplicate shared-string entry.
$ws.Range("B2").Value = "dương tính"
